$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row becomes text labels "X" / "Z" ---
$ws.Range("A1").Value = "X"
$ws.Range("B1").Value = "Z"

# --- Reverse the upper-surface block (rows 2:66) in place: ---
# Currently ascending by column A (trailing edge -> leading edge),
# re-sort descending so it now reads leading edge -> trailing edge.
$sortUpper = $ws.Sort
$sortUpper.SortFields.Clear()
$sortUpper.SortFields.Add($ws.Range("A2:A66"), 0, 2)
$sortUpper.SetRange($ws.Range("A2:B66"))
$sortUpper.Header = -4142
$sortUpper.Apply()

# --- Swap the two leading-edge duplicate rows (67 <-> 132) that sit ---
# --- just outside the lower-surface sort range. ---
$row67 = $ws.Range("A67:B67").Value2
$row132 = $ws.Range("A132:B132").Value2
$ws.Range("A67:B67").Value = $row132
$ws.Range("A132:B132").Value = $row67

# --- Reverse the lower-surface block (rows 68:131) in place: ---
# Currently descending by column A, re-sort ascending so it now reads
# leading edge -> trailing edge. This is the last sort applied, so it
# is the one captured by the sheet's persisted sortState.
$sortLower = $ws.Sort
$sortLower.SortFields.Clear()
$sortLower.SortFields.Add($ws.Range("A68:A131"), 0, 1)
$sortLower.SetRange($ws.Range("A68:B131"))
$sortLower.Header = -4142
$sortLower.Apply()

# --- Selection / view bookkeeping to match the saved file ---
[void]$ws.Range("H74").Select()
